# HelloAngular.docx edit: add more content ("libs") at the end of the
# document (new paragraphs about jQuery vs AngularJS two-way data binding,
# a code sample, and a closing reference line), per commit
# "add more libs and done for chapter10".

$d = $word.ActiveDocument

# The _GoBack bookmark sits (zero-length) in the last paragraph of the
# document.  Remove it now so our text insertions aren't perturbed by it;
# we re-create it later, right after its proper anchor text, once the new
# content has been typed in.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------
# 1) Fill in the previously-empty last paragraph (style "a4", ind=720)
#    with the "if you used jQuery..." sentence.
# ---------------------------------------------------------------------
$pCount = $d.Paragraphs.Count
$p = $d.Paragraphs($pCount)
$r = $p.Range
$r.Collapse(0)
$r.InsertBefore("如果使用")

$p = $d.Paragraphs($pCount)
$r = $p.Range
$r.Collapse(0)
$r.InsertBefore("jQuery")

$p = $d.Paragraphs($pCount)
$r = $p.Range
$r.Collapse(0)
$r.InsertBefore("我们需要给键盘绑定监控事件，当键盘被松开时读取输入框的内容同时显示在对应的DIV中。是不是想想都觉得很烦 ")

# Closing smiley -- real Word would store this as <w:sym w:font="Wingdings"
# w:char="F04A"/>; this runtime has no InsertSymbol support, so approximate
# it with the equivalent Private Use Area glyph set in the Wingdings font.
$p = $d.Paragraphs($pCount)
$r = $p.Range
$r.Collapse(0)
$r.InsertBefore([string][char]0xF04A)
$p = $d.Paragraphs($pCount)
$symRange = $d.Range($p.Range.End - 2, $p.Range.End - 1)
$symRange.Font.Name = "Arial Unicode MS"

# ---------------------------------------------------------------------
# Helper pattern used repeatedly below: append a brand new paragraph
# (inherits the style/indent of the paragraph before it), then type its
# text.
# ---------------------------------------------------------------------

function New-Para {
    $n = $d.Paragraphs.Count
    $p = $d.Paragraphs($n)
    $r = $p.Range
    $r.Collapse(0)
    $r.InsertParagraphAfter()
}

# 2) "但是在angular Js的世界里，这个行为被优雅的解决了。"
New-Para
$n = $d.Paragraphs.Count
$p = $d.Paragraphs($n); $r = $p.Range; $r.Collapse(0)
$r.InsertBefore("但是在angular ")
$p = $d.Paragraphs($n); $r = $p.Range; $r.Collapse(0)
$r.InsertBefore("Js")
$p = $d.Paragraphs($n); $r = $p.Range; $r.Collapse(0)
$r.InsertBefore("的世界里，这个行为被优雅的解决了。")

# 3) "我们只需要利用到angular Js提供的一个非常赞的特性“双向数据绑定”既可以完成这一需求。"
New-Para
$n = $d.Paragraphs.Count
$p = $d.Paragraphs($n); $r = $p.Range; $r.Collapse(0)
$r.InsertBefore("我们只需要利用到angular ")
$p = $d.Paragraphs($n); $r = $p.Range; $r.Collapse(0)
$r.InsertBefore("Js")
$p = $d.Paragraphs($n); $r = $p.Range; $r.Collapse(0)
$r.InsertBefore("提供的一个")
$p = $d.Paragraphs($n); $r = $p.Range; $r.Collapse(0)
$r.InsertBefore("非常")
$p = $d.Paragraphs($n); $r = $p.Range; $r.Collapse(0)
$r.InsertBefore("赞的特性“双向数据绑定”既可以完成这一需求。")

# 4) "代码片段："
New-Para
$n = $d.Paragraphs.Count
$p = $d.Paragraphs($n); $r = $p.Range; $r.Collapse(0)
$r.InsertBefore("代码片段：")

# 5) code line: <input type="text" ng-model="world" placeholder="input some words"/>
New-Para
$n = $d.Paragraphs.Count
$p = $d.Paragraphs($n); $r = $p.Range; $r.Collapse(0)
$r.InsertBefore('<input type="text" ')
$p = $d.Paragraphs($n); $r = $p.Range; $r.Collapse(0)
$r.InsertBefore("ng")
$p = $d.Paragraphs($n); $r = $p.Range; $r.Collapse(0)
$r.InsertBefore('-model="world"')
$p = $d.Paragraphs($n); $r = $p.Range; $r.Collapse(0)
$r.InsertBefore(' placeholder="input some words"/>')

# 6) code line: Hello {{world}}
New-Para
$n = $d.Paragraphs.Count
$p = $d.Paragraphs($n); $r = $p.Range; $r.Collapse(0)
$r.InsertBefore("Hello {{")
$p = $d.Paragraphs($n); $r = $p.Range; $r.Collapse(0)
$r.InsertBefore("world")
$p = $d.Paragraphs($n); $r = $p.Range; $r.Collapse(0)
$r.InsertBefore("}}")

# 7) "通过给input输入框绑定ng-model(angular js提供的一个内置指令)，..."
New-Para
$n = $d.Paragraphs.Count
$p = $d.Paragraphs($n); $r = $p.Range; $r.Collapse(0)
$r.InsertBefore("通过给input输入框绑定")
$p = $d.Paragraphs($n); $r = $p.Range; $r.Collapse(0)
$r.InsertBefore("ng")
$p = $d.Paragraphs($n); $r = $p.Range; $r.Collapse(0)
$r.InsertBefore("-model(angular ")
$p = $d.Paragraphs($n); $r = $p.Range; $r.Collapse(0)
$r.InsertBefore("js")
$p = $d.Paragraphs($n); $r = $p.Range; $r.Collapse(0)
$r.InsertBefore("提供的一个内置指令)，同时在页面上定义{{model}}，在两个花括号中间申明这个模型的值，这样就完成了数据的双向绑定。你在input内输入的任何内容都会在页面上即时的显示出来。")

# 8) "这样的双向数据绑定意味着，...剩下的时交给angularjs就可以了。" + re-homed _GoBack bookmark
New-Para
$n = $d.Paragraphs.Count
$p = $d.Paragraphs($n); $r = $p.Range; $r.Collapse(0)
$r.InsertBefore("这样的双向数据绑定意味着，你不在需要知道如何从模型上取值，更新值的内容后再重新的刷新HTML页面。而是，你仅仅只需要绑定他们，剩下的时交给")
$p = $d.Paragraphs($n); $r = $p.Range; $r.Collapse(0)
$r.InsertBefore("angularjs")
$p = $d.Paragraphs($n); $r = $p.Range; $r.Collapse(0)
$r.InsertBefore("就可以了。")

$p = $d.Paragraphs($n)
$bkPoint = $d.Range($p.Range.End - 1, $p.Range.End - 1)
$d.Bookmarks.Add("_GoBack", $bkPoint)

# 9) "详细代码参加HelloWorld.html"
New-Para
$n = $d.Paragraphs.Count
$p = $d.Paragraphs($n); $r = $p.Range; $r.Collapse(0)
$r.InsertBefore("详细代码参加HelloWorld.")
$p = $d.Paragraphs($n); $r = $p.Range; $r.Collapse(0)
$r.InsertBefore("html")

# ---------------------------------------------------------------------
# Formatting pass: red font color on the angular-binding placeholders in
# the code sample, matching the source's highlighted "ng-model"/"world".
# ---------------------------------------------------------------------

$codeLinePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs($i).Range.Text
    if ($txt -like '*ng-model=`"world`"*') {
        $codeLinePara = $i
    }
}
if ($codeLinePara -ne $null) {
    $scope = $d.Paragraphs($codeLinePara).Range
    $scope.Find.ClearFormatting()
    if ($scope.Find.Execute("ng-model=")) {
        $ngEnd = $scope.End
        $ngStart = $ngEnd - 2
        $d.Range($ngStart, $ngEnd).Font.Color = 255
        $modelStart = $ngEnd
        $modelEnd = $modelStart + 13
        $d.Range($modelStart, $modelEnd).Font.Color = 255
    }
}

$helloLinePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs($i).Range.Text
    if ($txt -like 'Hello {{world}}*') {
        $helloLinePara = $i
    }
}
if ($helloLinePara -ne $null) {
    $scope = $d.Paragraphs($helloLinePara).Range
    $scope.Find.ClearFormatting()
    if ($scope.Find.Execute("world")) {
        $d.Range($scope.Start, $scope.End).Font.Color = 255
    }
}
